# Update gh-pages to output generated at 456a3b4
# F3: 273 -> 276, F6: 788 -> 792 on the "展览" and "全部类型" sheets.
$wb = $excel.ActiveWorkbook

$targetSheetNames = @("展览", "全部类型")

foreach ($ws in $wb.Worksheets) {
    if ($targetSheetNames -contains $ws.Name) {
        $ws.Range("F3").Value = 276
        $ws.Range("F6").Value = 792
    }
}
